$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Recommended Readings" column (F) and the one extra cell in G2
# with corrected chapter/page references (learning objectives update).
$ws.Range("F2").Value  = "BMLR Ch. 2 (pp. 39-68), 6.1-6.4 (pp. 151-155)"
$ws.Range("G2").Value  = "BMLR 3.3 'Discrete Random Variables' (pp. 72-79)"
$ws.Range("F4").Value  = "BMLR Ch. 6.5-6.7 (pp. 159-180)"
$ws.Range("F8").Value  = "BMLR Ch. 4.1-4.4 (pp. 93-112)"
$ws.Range("F10").Value = "BMLR Ch. 4.5-4.10 (pp. 113-132), Ch. 5 (pp. 145-148)"
$ws.Range("F12").Value = "BMLR Ch. 8.1-8.5 (pp. 211-231)"
$ws.Range("F14").Value = "BMLR Ch. 8.6-8.11 (pp. 234-251)"
$ws.Range("F16").Value = "BMLR Ch. 9.1-9.7 (pp. 263-306), Ch. 11 (pp. 373-398)"
$ws.Range("F18").Value = "BMLR 7.1-7.8 (pp. 193-206)"
$ws.Range("F19").Value = "IEPM Ch. 5 & 6"

# Restore the view: zoom in to 160%, scroll so row 3 is near the top,
# and leave the selection on F19.
$win = $excel.ActiveWindow
$ws.Range("F19").Select()
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Zoom = 160
